$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 ("מיקום תוכנה" / software location) now points at the git working copy
# path instead of the old local user path (the old duplicate value that used
# to live in B5).
$ws.Range("B2").Value = "W:\Data\Forecast\Tools\forecast_git\create_forecast_basic\current"

# B5 held that duplicate path text - clear it out completely (value + format)
# now that its content has been promoted into B2.
$ws.Range("B5").Clear()

# The sheet had a stray far-away blank cell at B10 with several empty rows in
# between; tidy that up by deleting the empty filler rows 6:9, which pulls
# the trailing blank cell up from B10 to B6 and shrinks the used range.
$ws.Range("A6:B9").EntireRow.Delete()

# Reflect the new selection the user ended up with.
$ws.Range("B2:B3").Select()
